$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "'94.303.82"
$ws.Cells.Item(2, 4).Style = 'Normal'
$ws.Cells.Item(2, 5).Value = '  +2.75%  '
$ws.Cells.Item(3, 4).Value = "'3.085.61"
$ws.Cells.Item(3, 4).Style = 'Normal'
$ws.Cells.Item(3, 5).Value = '  +0.27%  '
$ws.Cells.Item(4, 5).Value = '  -0.08%  '
$ws.Cells.Item(5, 4).Value = "'237.05"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value = '  +1.10%  '
$ws.Cells.Item(6, 4).Value = "'610.56"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value = '  +0.48%  '
$ws.Cells.Item(7, 5).Value = '  +3.28%  '
$ws.Cells.Item(8, 4).Value = "'0.380"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value = '  -0.48%  '
$ws.Cells.Item(9, 5).Value = '  -0.06%  '
$ws.Cells.Item(10, 4).Value = "'0.813"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value = '  +12.46%  '
$ws.Cells.Item(11, 4).Value = "'3.081.95"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value = '  +0.34%  '
$ws.Cells.Item(12, 5).Value = '  -1.41%  '
$ws.Cells.Item(13, 4).Value = "'94.051.39"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value = '  +2.16%  '
$ws.Cells.Item(14, 4).Value = "'0.0000242"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value = '  -1.03%  '
$ws.Cells.Item(15, 4).Value = "'34.16"
$ws.Cells.Item(15, 4).Style = 'Normal'
$ws.Cells.Item(15, 5).Value = '  +1.53%  '
$ws.Cells.Item(16, 4).Value = "'5.33"
$ws.Cells.Item(16, 4).Style = 'Normal'
$ws.Cells.Item(16, 5).Value = '  -0.89%  '
$ws.Cells.Item(17, 4).Value = "'3.661.20"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value = '  -0.36%  '
$ws.Cells.Item(18, 4).Value = "'3.071.62"
$ws.Cells.Item(18, 4).Style = 'Normal'
$ws.Cells.Item(18, 5).Value = '  -0.75%  '
$ws.Cells.Item(19, 4).Value = "'3.58"
$ws.Cells.Item(19, 4).Style = 'Normal'
$ws.Cells.Item(19, 5).Value = '  -2.23%  '
$ws.Cells.Item(20, 4).Value = "'14.49"
$ws.Cells.Item(20, 4).Style = 'Normal'
$ws.Cells.Item(20, 5).Value = '  +0.01%  '
$ws.Cells.Item(21, 4).Value = "'5.78"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value = '  +1.50%  '
$ws.Cells.Item(22, 4).Value = "'448.41"
$ws.Cells.Item(22, 4).Style = 'Normal'
$ws.Cells.Item(22, 5).Value = '  +1.94%  '
$ws.Cells.Item(23, 4).Value = "'8.87"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value = '  -3.69%  '
$ws.Cells.Item(24, 5).Value = '  -0.47%  '
$ws.Cells.Item(25, 4).Value = "'8.34"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value = '  +6.42%  '
$ws.Cells.Item(26, 4).Value = "'5.54"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value = '  -1.84%  '
$ws.Cells.Item(27, 4).Value = "'12.02"
$ws.Cells.Item(27, 4).Style = 'Normal'
$ws.Cells.Item(27, 5).Value = '  +4.86%  '
$ws.Cells.Item(28, 4).Value = "'84.93"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value = '  -0.82%  '
$ws.Cells.Item(29, 4).Value = "'3.244.35"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value = '  -0.41%  '
$ws.Cells.Item(30, 4).Value = "'1.00"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value = '  +0.24%  '
$ws.Cells.Item(31, 5).Value = '  +11.76%  '
$ws.Cells.Item(32, 4).Value = "'0.179"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value = '  +7.42%  '
$ws.Cells.Item(33, 5).Value = '  -4.49%  '
$ws.Cells.Item(34, 5).Value = '  +0.37%  '
$ws.Cells.Item(35, 4).Value = "'0.998"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value = '  +0.54%  '
$ws.Cells.Item(36, 4).Value = "'7.67"
$ws.Cells.Item(36, 4).Style = 'Normal'
$ws.Cells.Item(36, 5).Value = '  -0.82%  '
$ws.Cells.Item(37, 2).Value = 'EthereumClassic'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(37, 4).Value = "'25.55"
$ws.Cells.Item(37, 4).Style = 'Normal'
$ws.Cells.Item(37, 5).Value = '  -0.34%  '
$ws.Cells.Item(38, 2).Value = 'Kaspa'
$ws.Cells.Item(38, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(38, 4).Value = "'0.152"
$ws.Cells.Item(38, 4).Style = 'Normal'
$ws.Cells.Item(38, 5).Value = '  -2.67%  '
$ws.Cells.Item(39, 5).Value = '  +0.86%  '
$ws.Cells.Item(40, 4).Value = "'481.95"
$ws.Cells.Item(40, 4).Style = 'Normal'
$ws.Cells.Item(40, 5).Value = '  +1.11%  '
$ws.Cells.Item(41, 4).Value = "'24.05"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value = '  +0.87%  '
$ws.Cells.Item(42, 4).Value = "'0.439"
$ws.Cells.Item(42, 4).Style = 'Normal'
$ws.Cells.Item(42, 5).Value = '  +3.12%  '
$ws.Cells.Item(43, 2).Value = 'MantraDAO'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Cells.Item(43, 4).Value = "'3.71"
$ws.Cells.Item(43, 4).Style = 'Normal'
$ws.Cells.Item(43, 5).Value = '  -4.07%  '
$ws.Cells.Item(44, 2).Value = 'Fetch.AI'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(44, 4).Value = "'1.25"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value = '  -1.71%  '
$ws.Cells.Item(46, 4).Value = "'3.10"
$ws.Cells.Item(46, 4).Style = 'Normal'
$ws.Cells.Item(46, 5).Value = '  -4.18%  '
$ws.Cells.Item(47, 4).Value = "'161.36"
$ws.Cells.Item(47, 4).Style = 'Normal'
$ws.Cells.Item(47, 5).Value = '  -0.83%  '
$ws.Cells.Item(48, 4).Value = "'0.676"
$ws.Cells.Item(48, 4).Style = 'Normal'
$ws.Cells.Item(48, 5).Value = '  -0.44%  '
$ws.Cells.Item(49, 4).Value = "'1.82"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value = '  -1.53%  '
$ws.Cells.Item(50, 4).Value = "'43.68"
$ws.Cells.Item(50, 4).Style = 'Normal'
$ws.Cells.Item(50, 5).Value = '  -0.50%  '
$ws.Cells.Item(51, 4).Value = "'0.999"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value = '  +0.17%  '
